$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSAD")
$ws.Range("A1").Value = "Hello"
$ws.Range("A2").Value = "World"
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:B2"), $null, 1)
$lo.Name = "TestTable"
